$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MCT-1A-Desenho Técnico"
$ws.Range("B3").Value = "[-, 'MCT-2A-CAD']"
$ws.Range("D3").Value = "MCT-3A-Máquinas Térmicas e de Fluxo"
$ws.Range("B4").Value = "[-, 'MCT-2A-CAD']"
$ws.Range("D4").Value = "MCT-3A-Máquinas Térmicas e de Fluxo"
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "-"
$ws.Range("E6").Value = "-"
$ws.Range("C7").Value = "-"
$ws.Range("E7").Value = "-"
$ws.Range("C8").Value = "MCT-1A-Desenho Técnico"
